$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match results for Round 25, appended as rows 131-139 right after the
# existing data (which ended at row 130). Each tuple is:
#   matchIndex, homeTeam, awayTeam, homeScore, awayScore, homeXg, awayXg,
#   homeXgPred, awayXgPred, homePred, awayPred, homeDiff, awayDiff, totalDiff, goalDiff
$data = @(
    @(129, "Cagliari", "Hellas Verona", 1, 0, 1.91, 0.7,  2.27, 0.74, 1, 0, 0.36, 0.04, 0.41, 0),
    @(130, "Bologna",  "Venezia",       3, 0, 3.22, 0.44, 3.37, 0.53, 2, 0, 0.15, 0.09, 0.24, 1),
    @(131, "Como",     "Monza",         1, 1, 1.29, 1.34, 1.8,  1.02, 0, 1, 0.51, 0.32, 0.82, 1),
    @(132, "Milan",    "Empoli",        3, 0, 2.29, 0.54, 2.24, 0.65, 0, 0, 0.05, 0.11, 0.16, 3),
    @(133, "Lecce",    "Juventus",      1, 1, 1.51, 1.3,  1.73, 1.07, 1, 1, 0.22, 0.23, 0.44, 0),
    @(134, "Parma",    "Lazio",         3, 1, 1.48, 2.6,  1.11, 2.82, 0, 1, 0.37, 0.22, 0.59, 3),
    @(135, "Torino",   "Napoli",        0, 1, 0.26, 1.41, 0.42, 1.74, 0, 0, 0.16, 0.33, 0.49, 1),
    @(136, "Udinese",  "Genoa",         0, 2, 0.17, 1.64, 0.25, 1.87, 0, 1, 0.08, 0.23, 0.31, 1),
    @(137, "Roma",     "Atalanta",      0, 2, 0.42, 1.3,  0.36, 1.37, 0, 0, 0.06, 0.07000000000000001, 0.14, 2)
)

$startRow = 131
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = $row[13]
    $ws.Cells.Item($r, 15).Value = $row[14]
}

# Column A uses the same bold/centered/bordered style as the rest of column A
# (e.g. A130). Copy that formatting down onto the newly added cells without
# touching the values we just wrote.
$ws.Range("A130").Copy()
$ws.Range("A131:A139").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
